$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.240.06"
$ws.Range('E2').Value = '  +1.02%  '

$ws.Range('D3').Value = "'1.802.23"
$ws.Range('E3').Value = '  +2.50%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = "'339.32"
$ws.Range('E5').Value = '  +0.21%  '

$ws.Range('D6').Value = "'0.9992"
$ws.Range('E6').Value = '  -0.37%  '

$ws.Range('D7').Value = "'0.4647"
$ws.Range('E7').Value = '  +23.76%  '

$ws.Range('D8').Value = "'0.3635"
$ws.Range('E8').Value = '  +8.45%  '

$ws.Range('D9').Value = "'45.55"
$ws.Range('E9').Value = '  -0.84%  '

$ws.Range('D10').Value = "'1.145"
$ws.Range('E10').Value = '  +2.13%  '

$ws.Range('D11').Value = "'0.07572"
$ws.Range('E11').Value = '  +6.40%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'22.52"
$ws.Range('E12').Value = '  +1.37%  '

$ws.Range('B13').Value = 'BinanceUSD'
$ws.Range('C13').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D13').Value = "'1.002"
$ws.Range('E13').Value = '  -0.53%  '

$ws.Range('D14').Value = "'6.250"
$ws.Range('E14').Value = '  +1.36%  '

$ws.Range('D15').Value = "'7.275"
$ws.Range('E15').Value = '  +1.81%  '

$ws.Range('D16').Value = "'1.793.29"
$ws.Range('E16').Value = '  +2.23%  '

$ws.Range('D17').Value = "'0.00001087"
$ws.Range('E17').Value = '  +3.86%  '

$ws.Range('D18').Value = "'0.06713"
$ws.Range('E18').Value = '  +2.09%  '

$ws.Range('D19').Value = "'81.51"
$ws.Range('E19').Value = '  +1.83%  '

$ws.Range('E20').Value = '  -0.24%  '

$ws.Range('D21').Value = "'17.27"
$ws.Range('E21').Value = '  +2.31%  '

$ws.Range('D22').Value = "'6.381"
$ws.Range('E22').Value = '  +2.12%  '

$ws.Range('D23').Value = "'28.220.39"
$ws.Range('E23').Value = '  +1.08%  '

$ws.Range('D24').Value = "'11.91"
$ws.Range('E24').Value = '  +2.46%  '

$ws.Range('D25').Value = "'2.405"
$ws.Range('E25').Value = '  +0.38%  '

$ws.Range('D26').Value = "'20.43"
$ws.Range('E26').Value = '  +4.04%  '

$ws.Range('D27').Value = "'2.405"
$ws.Range('E27').Value = '  +4.19%  '

$ws.Range('D28').Value = "'152.44"
$ws.Range('E28').Value = '  +0.24%  '

$ws.Range('D29').Value = "'1.998.87"
$ws.Range('E29').Value = '  +1.90%  '

$ws.Range('D30').Value = "'1.276"
$ws.Range('E30').Value = '  -0.14%  '

$ws.Range('D31').Value = "'132.82"
$ws.Range('E31').Value = '  +1.25%  '

$ws.Range('D32').Value = "'4.074"
$ws.Range('E32').Value = '  +1.15%  '

$ws.Range('D33').Value = "'5.906"
$ws.Range('E33').Value = '  +2.86%  '

$ws.Range('D34').Value = "'0.09517"
$ws.Range('E34').Value = '  +9.88%  '

$ws.Range('D35').Value = "'0.02374"
$ws.Range('E35').Value = '  +1.92%  '

$ws.Range('D36').Value = "'12.12"
$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('D37').Value = "'0.06280"
$ws.Range('E37').Value = '  +1.72%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'5.199"
$ws.Range('E38').Value = '  +1.84%  '

$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = "'0.6620"
$ws.Range('E39').Value = '  +1.68%  '

$ws.Range('D40').Value = "'0.2168"
$ws.Range('E40').Value = '  +3.30%  '

$ws.Range('D41').Value = "'1.480"
$ws.Range('E41').Value = '  +1.93%  '

$ws.Range('E42').Value = '  +0.24%  '

$ws.Range('D43').Value = "'8.081"
$ws.Range('E43').Value = '  +0.64%  '

$ws.Range('D44').Value = "'0.9989"
$ws.Range('E44').Value = '  -0.28%  '

$ws.Range('D45').Value = "'13.96"

$ws.Range('D46').Value = "'3.870"
$ws.Range('E46').Value = '  +0.90%  '

$ws.Range('D47').Value = "'0.6093"
$ws.Range('E47').Value = '  +1.82%  '

$ws.Range('D48').Value = "'128.56"
$ws.Range('E48').Value = '  -0.64%  '

$ws.Range('D49').Value = "'2.031"
$ws.Range('E49').Value = '  +1.80%  '

$ws.Range('D50').Value = "'0.07097"
$ws.Range('E50').Value = '  -1.44%  '

$ws.Range('D51').Value = "'1.169"
$ws.Range('E51').Value = '  +0.00%  '
